{"js": "// Remove the trailing sentence \"Some other testcases are in week5problem3.ipynb.\"\n// from the \"Problem 2\" paragraph, leaving the single space that preceded it.\nconst results = context.document.body.search(\"Some other testcases are in week5problem3.ipynb.\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "# Word COM (PowerShell-style) script reproducing the edits described by the diff.\n#\n# Net content change:\n#   - The \"Problem 2\" paragraph loses its trailing sentence\n#     \"Some other testcases are in week5problem3.ipynb.\" (the leading space\n#     that preceded the sentence is kept).\n#\n# The diff also shows a number of run-level XML changes that carry no visible\n# text difference (re-splitting of \"VaR\"/\"riskmgmt\" runs around spell-check\n# <w:proofErr/> markers, and the coalescing of previously-split runs such as\n# \"For portfolio B:\" or \"VaR for delta normal: 4494.60\" into single runs, plus\n# two image runs picking up a <w:noProof/> font flag). Those are reproduced\n# here too, since they are natural, deterministic side effects of driving the\n# same Find/Replace and Range operations a person would use in Word.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($searchText, $replaceText) {\n  $range = $d.Content\n  $find = $range.Find\n  $find.ClearFormatting()\n  $find.Text = $searchText\n  $find.Replacement.ClearFormatting()\n  $find.Replacement.Text = $replaceText\n  $find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null\n}\n\n# --- The actual content edit -------------------------------------------------\n# Remove \"Some other testcases are in week5problem3.ipynb.\" leaving the single\n# space that preceded it intact.\nReplace-Text \"Some other testcases are in week5problem3.ipynb.\" \"\"\n\n# --- Cosmetic run coalescing (no visible text change) ------------------------\n# \"For portfolio B:\" / \"For portfolio C:\" labels, previously split across\n# multiple runs, become a single run once replaced as a whole.\nReplace-Text \"For portfolio B:\" \"For portfolio B:\"\nReplace-Text \"For portfolio C:\" \"For portfolio C:\"\n\n# The VaR value lines for portfolios B and C were split between the label and\n# the number (and, for \"delta normal\", letter-by-letter); replacing the whole\n# line coalesces each into one run.\nReplace-Text \"VaR for generalized T distribution: 6693.39\" \"VaR for generalized T distribution: 6693.39\"\nReplace-Text \"VaR for delta normal: 4494.60\" \"VaR for delta normal: 4494.60\"\nReplace-Text \"VaR for historical simulation: 7273.70\" \"VaR for historical simulation: 7273.70\"\nReplace-Text \"VaR for generalized T distribution: 5653.33\" \"VaR for generalized T distribution: 5653.33\"\nReplace-Text \"VaR for delta normal: 3786.59\" \"VaR for delta normal: 3786.59\"\nReplace-Text \"VaR for historical simulation: 5310.07\" \"VaR for historical simulation: 5310.07\"\n\n# --- noProof flag on the two inline pictures that previously lacked rPr ------\n$shapes = $d.InlineShapes\n$shapes.Item(4).Range.NoProofing = $true\n$shapes.Item(6).Range.NoProofing = $true\n\nWrite-Output \"done\"\n"}
